# TC31_Canine_Filter_Breed-IrishSettr.xlsx
# Fixed variables and query errors in Breed from TC30 to TC47
#
# The "CasesTab" Cypher query (cell B2 on the "startup" sheet) had a stray
# trailing `co:cohort` / `Cohort` column that doesn't belong on this query
# (it duplicated logic that lives on a different tab's query). Remove it so
# the query ends cleanly after "Response to Treatment".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesTabQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" +
    "WHERE demo.breed IN ['Irish Setter']`n" +
    "MATCH (c)<--(diag:diagnosis)`n" +
    "OPTIONAL MATCH (samp:sample)-->(c)`n" +
    "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" +
    "WITH DISTINCT c, s, demo, diag, co`n" +
    "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" +
    "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" +
    "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" +
    "        coalesce(demo.breed, '') AS Breed ,`n" +
    "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" +
    "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" +
    "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" +
    "        coalesce(demo.sex, '') AS Sex ,`n" +
    "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" +
    "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" +
    "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value2 = $casesTabQuery

# The query text got shorter, so the wrapped row heights for the three data
# rows (CasesTab / SamplesTab / FilesTab) shrink accordingly.
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 216

# Selection moves to B2 (the edited cell) instead of B4, and the view
# scrolls back so row 1 is in view again.
$ws.Range("B2").Select()
